$d = $word.ActiveDocument

$wdYellow = 7

function Highlight-Text($text) {
    $r = $d.Content
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $text"
        return
    }
    $r.Font.HighlightColorIndex = $wdYellow
    Write-Output "Highlighted: [$text]"
}

# Paragraph: "Ajout d'un compte administrateur. Vous pouvez ... administrateur. "
# -> split into two highlighted sentence runs, separated/trailed by plain spaces.
$t1 = "Ajout d" + [char]0x2019 + "un compte administrateur."
Highlight-Text $t1

$t2 = "Vous pouvez mettre cette information directement dans la base de donn" + [char]0x00E9 + "es; vous n" + [char]0x2019 + "avez pas " + [char]0x00E0 + " permettre la cr" + [char]0x00E9 + "ation de comptes de type administrateur."
Highlight-Text $t2

# Paragraph: "Les routes liées aux comptes auront cette forme : /comptes/XXXX "
$t3 = "Les routes li" + [char]0x00E9 + "es aux comptes auront cette forme : /comptes/XXXX"
Highlight-Text $t3

# Paragraph: "Un utilisateur a un courriel et un mot de passe. ... minimum. "
$t4 = "Un utilisateur a un courriel et un mot de passe. Le courriel doit " + [char]0x00EA + "tre valide et le mot de passe doit respecter les r" + [char]0x00E8 + "gles suivantes : Une lettre majuscule, une lettre minuscule, un nombre et avoir une longueur de 8 charact" + [char]0x00E8 + "res au minimum."
Highlight-Text $t4

Write-Output "done"
